$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Julia", "Álvarez", 5161149, 3022154659, "julia.alvarez@gmail.com"),
    @("Luisa", "Hernández", 4659873, 3084531624, "luisa.hernanza@hotmail.com"),
    @("Mauricio", "Rodríguez", 2152426, 3105498675, "camilo.rodri@gmail.com"),
    @("Pablo", "Casas", 2356849, 3152468975, "p.casas@gmail.com"),
    @("Ángela", "Ruiz", 2114853, 3002586491, "angela.r@gmail.com")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
}

$ws.Columns.Item(4).ColumnWidth = 13.17
